$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow programmatic edits, then
# restore protection with the same password once done.
$ws.Unprotect("D382")

# Update the "as of" date in the confidentiality / disclosure footer text.
$ws.Cells.Item(42, 1).Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-07 for illustrative purposes only and are subject to change."

# Refresh the Weight (column D) and Percent Change (column E) values for
# each holding row (rows 2-38) and the Percent Change total (row 39).
$ws.Cells.Item(2, 4).Value = 0.05891157382083293
$ws.Cells.Item(2, 5).Value = 0.003622629875134953
$ws.Cells.Item(3, 4).Value = 0.05274229077215829
$ws.Cells.Item(3, 5).Value = 0.010931806350859
$ws.Cells.Item(4, 4).Value = 0.3079467191562247
$ws.Cells.Item(4, 5).Value = 0.02352376380220833
$ws.Cells.Item(5, 4).Value = 0.03491481358674188
$ws.Cells.Item(5, 5).Value = -0.004464110187305081
$ws.Cells.Item(6, 4).Value = 0.0317551970364801
$ws.Cells.Item(6, 5).Value = 0.003458412588621895
$ws.Cells.Item(7, 4).Value = 0.03054355837808955
$ws.Cells.Item(7, 5).Value = 0.003422739436181432
$ws.Cells.Item(8, 4).Value = 0.02834098219395934
$ws.Cells.Item(8, 5).Value = 0.004530821509478811
$ws.Cells.Item(9, 4).Value = 0.02383149838117303
$ws.Cells.Item(9, 5).Value = -0.006026231832683582
$ws.Cells.Item(10, 4).Value = 0.02468209533021747
$ws.Cells.Item(10, 5).Value = 0.00623783344385731
$ws.Cells.Item(11, 4).Value = 0.0236555710668207
$ws.Cells.Item(11, 5).Value = -0.002937316417723879
$ws.Cells.Item(12, 4).Value = 0.02306823149753639
$ws.Cells.Item(12, 5).Value = 0.004046655558200385
$ws.Cells.Item(13, 4).Value = 0.0211164520550549
$ws.Cells.Item(13, 5).Value = 0.01677760052808197
$ws.Cells.Item(14, 4).Value = 0.02140589789457455
$ws.Cells.Item(14, 5).Value = 0.00007399733609592118
$ws.Cells.Item(15, 4).Value = 0.02138879091682841
$ws.Cells.Item(15, 5).Value = 0.004946975531725739
$ws.Cells.Item(16, 4).Value = 0.02191013134357359
$ws.Cells.Item(16, 5).Value = 0.00645829819023036
$ws.Cells.Item(17, 4).Value = 0.01982508643247707
$ws.Cells.Item(17, 5).Value = -0.0002130606157452197
$ws.Cells.Item(18, 4).Value = 0.01449404529336355
$ws.Cells.Item(18, 5).Value = 0.008393075712537224
$ws.Cells.Item(19, 4).Value = 0.01684762751562332
$ws.Cells.Item(19, 5).Value = 0.01983151983151998
$ws.Cells.Item(20, 4).Value = 0.01565235664458247
$ws.Cells.Item(20, 5).Value = -0.009613762860516073
$ws.Cells.Item(21, 4).Value = 0.01689894844886175
$ws.Cells.Item(21, 5).Value = 0.01429731925264011
$ws.Cells.Item(22, 4).Value = 0.01401378273293473
$ws.Cells.Item(22, 5).Value = 0.01330741176115979
$ws.Cells.Item(23, 4).Value = 0.01497430785379236
$ws.Cells.Item(23, 5).Value = -0.0005500550055005382
$ws.Cells.Item(24, 4).Value = 0.0147165416026915
$ws.Cells.Item(24, 5).Value = -0.007713668620796033
$ws.Cells.Item(25, 4).Value = 0.01384398013901003
$ws.Cells.Item(25, 5).Value = 0.01887871853546907
$ws.Cells.Item(26, 4).Value = 0.01379846713031504
$ws.Cells.Item(26, 5).Value = 0.00134691472346149
$ws.Cells.Item(27, 4).Value = 0.01268672477407174
$ws.Cells.Item(27, 5).Value = 0.03236197467975122
$ws.Cells.Item(28, 4).Value = 0.01381483491766475
$ws.Cells.Item(28, 5).Value = 0.009172628875435684
$ws.Cells.Item(29, 4).Value = 0.01407059479483243
$ws.Cells.Item(29, 5).Value = 0.009951518244450064
$ws.Cells.Item(30, 4).Value = 0.01347269536274185
$ws.Cells.Item(30, 5).Value = 0.008276899924755554
$ws.Cells.Item(31, 4).Value = 0.01229632664612421
$ws.Cells.Item(31, 5).Value = 0.00003435127614981326
$ws.Cells.Item(32, 4).Value = 0.01348367762006036
$ws.Cells.Item(32, 5).Value = -0.002842866988283843
$ws.Cells.Item(33, 4).Value = 0.01232019193606637
$ws.Cells.Item(33, 5).Value = 0.008099768578040401
$ws.Cells.Item(34, 4).Value = 0.006134435501413966
$ws.Cells.Item(34, 5).Value = 0.01991668388074097
$ws.Cells.Item(35, 4).Value = 0.005275179464868393
$ws.Cells.Item(35, 5).Value = 0.008587728956060303
$ws.Cells.Item(36, 4).Value = 0.005322593248868511
$ws.Cells.Item(36, 5).Value = 0.005317038330291179
$ws.Cells.Item(37, 4).Value = 0.00510685525173657
$ws.Cells.Item(37, 5).Value = 0.01058704327867499
$ws.Cells.Item(38, 4).Value = 0.004736943257633196
$ws.Cells.Item(38, 5).Value = 0.008917027063177096
$ws.Cells.Item(39, 5).Value = 0.01088848573679324

$ws.Protect("D382")
